# OLX Monitor run 2026-02-22 21:03 — append the latest scrape results to the
# "PODSUMOWANIE" sheet (the long, row-per-listing log that the diff touches:
# dimension A1:H194 -> A1:H202, 8 new rows appended after row 194).
#
# The 8 new rows are a verbatim re-post of the previous 8 rows (same
# listings/prices/urls as rows 187-194) with only the "last checked"
# timestamp in column A refreshed to the new run time — exactly what the
# monitoring script does when a listing is still live on the next check.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

$newTimestamp = "2026-02-22 21:03:42"
$templateRows = 8

# Find the current last used row (194 before this edit) so the template
# block and the destination are computed relative to the real sheet state.
$lastRow = $ws.UsedRange.Rows.Count
$firstTemplateRow = $lastRow - $templateRows + 1
$destFirstRow = $lastRow + 1
$destLastRow = $lastRow + $templateRows

# Copy the last 8 rows (values + styles) straight down, producing rows
# 195-202 with identical content/formatting to rows 187-194.
$srcRange = $ws.Range("A" + $firstTemplateRow + ":H" + $lastRow)
$destCell = $ws.Range("A" + $destFirstRow)
$srcRange.Copy($destCell)

# Stamp the new rows with this run's "last checked" timestamp.
$ws.Range("A" + $destFirstRow + ":A" + $destLastRow).Value = $newTimestamp
